# Add four new CPRA gate/station locations (rows 13-16) to the tracker sheet.
# Column A = short station code, Column B = friendly display name,
# Column C = numeric count, defaulted to 0 for brand-new stations.
#
# Column A is filled fully before column B so new shared-string entries are
# interned in the same order produced by the original authoring tool
# (all-of-A, then all-of-B, then C) - this keeps the shared-strings table
# layout identical to the target.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "LakefrontAirport"
$ws.Range("A14").Value = "Mandeville"
$ws.Range("A15").Value = "Rigolets"
$ws.Range("A16").Value = "Lafitte"

$ws.Range("B13").Value = "Lakefront Airport"
$ws.Range("B14").Value = "Mandeville"
$ws.Range("B15").Value = "Rigolets"
$ws.Range("B16").Value = "Lafitte"

$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("C13:C16").NumberFormat = "0.0"

# Widen column A to fit the new (longer) station codes, matching the
# bestFit width Excel would compute for "LakefrontAirport".
$ws.Columns("A:A").ColumnWidth = 14.6640625
